# Rename the inline picture shapes (the drawing's "name", i.e. the
# wp:docPr/@name and pic:cNvPr/@name attributes) for the three logo
# images that live in the document's header/footer stories:
#   - header (BTec_Logo-Orange)   : image1.jpg -> image2.jpg
#   - footer #1 (PearsonLogo, id=3): image2.png -> image1.png
#   - footer #2 (PearsonLogo, id=2): image2.png -> image1.png
#
# InlineShape has no settable "Name" that round-trips onto the picture's
# non-visual properties (pic:cNvPr), only onto wp:docPr, so we rebuild the
# <w:drawing> run via Range.InsertXML (preserving every other attribute
# verbatim) and then delete the original drawing run that InsertXML leaves
# in place ahead of the freshly inserted one.

$d = $word.ActiveDocument

function Rewrite-InlineLogo {
    param(
        [string]$MatchDescr,
        [string]$OldName,
        [string]$NewName,
        [string]$Id,
        [string]$Cx,
        [string]$Cy
    )

    foreach ($story in $d.StoryRanges) {
        if ($story.InlineShapes.Count -le 0) { continue }

        for ($i = 1; $i -le $story.InlineShapes.Count; $i++) {
            $shape = $story.InlineShapes($i)
            if ($shape.AlternativeText -ne $MatchDescr) { continue }

            $descr = $MatchDescr

            $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                '<pkg:xmlData>' +
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' +
                '<w:body><w:p><w:r><w:drawing>' +
                '<wp:inline distB="0" distT="0" distL="0" distR="0">' +
                '<wp:extent cx="' + $Cx + '" cy="' + $Cy + '"/>' +
                '<wp:effectExtent b="0" l="0" r="0" t="0"/>' +
                '<wp:docPr descr="' + $descr + '" id="' + $Id + '" name="' + $NewName + '"/>' +
                '<a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
                '<pic:pic><pic:nvPicPr>' +
                '<pic:cNvPr descr="' + $descr + '" id="0" name="' + $NewName + '"/>' +
                '<pic:cNvPicPr preferRelativeResize="0"/>' +
                '</pic:nvPicPr>' +
                '<pic:blipFill><a:blip r:embed="rId1"/><a:srcRect b="0" l="0" r="0" t="0"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill>' +
                '<pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="' + $Cx + '" cy="' + $Cy + '"/></a:xfrm><a:prstGeom prst="rect"/><a:ln/></pic:spPr>' +
                '</pic:pic></a:graphicData></a:graphic>' +
                '</wp:inline></w:drawing></w:r></w:p></w:body></w:document>' +
                '</pkg:xmlData></pkg:part></pkg:package>'

            $shape.Range.InsertXML($xml)

            # InsertXML dropped the freshly-built drawing right before the
            # original one; the original (still carrying $OldName) is now
            # the first inline shape in this story, so remove it.
            $story.InlineShapes(1).Delete()
        }
    }
}

# header1.xml: BTec_Logo-Orange, id=1, 914400 x 277792 EMU
Rewrite-InlineLogo -MatchDescr "BTec_Logo-Orange" -OldName "image1.jpg" -NewName "image2.jpg" -Id "1" -Cx "914400" -Cy "277792"

# footer1.xml: PearsonLogo, id=3, 952500 x 285750 EMU
Rewrite-InlineLogo -MatchDescr "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" -OldName "image2.png" -NewName "image1.png" -Id "3" -Cx "952500" -Cy "285750"

# footer2.xml: PearsonLogo, id=2, 952500 x 285750 EMU
Rewrite-InlineLogo -MatchDescr "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" -OldName "image2.png" -NewName "image1.png" -Id "2" -Cx "952500" -Cy "285750"

Write-Output "done"
